# Add a new worksheet "getAllCacheNames" with a connector-cache test case,
# matching the structure of the existing "getConceptModelDataByCondition" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet goes at the very end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "getAllCacheNames"

# Copy the header formatting (fill/border/font) from sheet1's header row so the
# new sheet's header row looks the same (gray header style) without touching values.
$ws1.Range("A1:F1").Copy() | Out-Null
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null

# Header row (reuses existing shared strings where the text already exists).
$newSheet.Range("A1").Value = "test-id"
$newSheet.Range("B1").Value = "description"
$newSheet.Range("D1").Value = "rspStatus"
$newSheet.Range("E1").Value = "rspCode"
$newSheet.Range("F1").Value = "rspMessage"

# Data row.
$newSheet.Range("A2").Value = "jinzu-connector-cache-test1"
$newSheet.Range("B2").Value = "get allCacheNames"

# New header text (written after the data row so the shared-string table is
# appended in the same order the source workbook used).
$newSheet.Range("C1").Value = "dataList"

$newSheet.Range("C2").Value = "table_schema,datasource_config,datasource_restTemplate,datasource_connection,mapper,datasource_auth,datasource_freemarker,datasource_webservice,executorBaseInfo,plugin_relatedDb,plugin_thriftClient,plugin_debug,entityCache,data,meta,rule"

# Column widths (closest values reachable through the ColumnWidth property).
$newSheet.Columns.Item(1).ColumnWidth = 27
$newSheet.Columns.Item(2).ColumnWidth = 25.5714285714
$newSheet.Columns.Item(3).ColumnWidth = 25.5714285714

# sheet1 keeps a scrolled-over selection, but loses the active-tab flag -
# set this first since selecting on a sheet activates it.
$ws1.Activate() | Out-Null
$ws1.Range("L2").Select() | Out-Null

# Selection on the new sheet, then make it the active (last-selected) tab.
$newSheet.Range("C11").Select() | Out-Null
$newSheet.Activate() | Out-Null
